$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 (int / String) and shift the rows below it up
$ws.Rows.Item(4).Delete()

# Update the selection to match the saved state (B4)
$ws.Range("B4").Select()
